$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.209.00'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -0.78%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.172.73'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -4.52%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '591.71'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -2.24%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.91'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -5.40%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.170.01'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -4.54%  '
$ws.Range("E9").Value = '  -0.85%  '
$ws.Range("E10").Value = '  -5.63%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.25'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -5.09%  '
$ws.Range("E12").Value = '  -3.05%  '
$ws.Range("E13").Value = '  -4.29%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.00'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +0.15%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.695.19'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -4.49%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.172.61'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -4.36%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '63.144.33'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -1.01%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.60'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -3.95%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '461.39'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -3.97%  '
$ws.Range("E21").Value = '  -1.37%  '
$ws.Range("E22").Value = '  -5.02%  '
$ws.Range("E23").Value = '  -4.22%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.51'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -1.73%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.44'
$ws.Range("D25").ClearFormats()
$ws.Range("E26").Value = '  -0.02%  '
$ws.Range("E27").Value = '  +0.12%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.68'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -3.76%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.75'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -6.13%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.77'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -5.86%  '
$ws.Range("E31").Value = '  -5.84%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '27.21'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -6.09%  '
$ws.Range("E33").Value = '  -4.12%  '
$ws.Range("E34").Value = '  -6.40%  '
$ws.Range("E35").Value = '  -6.09%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.85'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -3.99%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '51.42'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -1.86%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0₃0708'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -4.92%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0389'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -2.91%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '405.65'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -6.79%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.13'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -2.43%  '
$ws.Range("E42").Value = '  -5.05%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.817.76'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -9.10%  '
$ws.Range("E44").Value = '  -6.26%  '
$ws.Range("E45").Value = '  -5.74%  '
$ws.Range("E46").Value = '  -5.62%  '
$ws.Range("E47").Value = '  -0.06%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '25.39'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -3.98%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '124.04'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -0.27%  '
$ws.Range("E50").Value = '  -1.81%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '34.00'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -8.37%  '
